$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# The sheet lists medicines alphabetically (row 4 .. row 23 originally).
# Three new medicines were added:
#   BETOLVEX 1MG/ML 2 AMP           -> balance 1:1, price 64, txns 1:0
#   EPIMAG EFFERVESCENT 12 SACHETS  -> balance 0:0, price 26, txns 1:0
#   HERO VITAMIN DROPS              -> balance 1:0, price 32, txns 1:0
# which pushes the alphabetically-later rows down. Row 4 (BEBY RELIEF)
# keeps its place; everything from row 5 on is rewritten in place (same
# row heights / styles as before) and three brand-new rows (24-26) are
# appended before the totals row.
# -----------------------------------------------------------------

# Final table, in row order, after the 3 new items are inserted.
$finalRows = @(
    ,@('BEBY RELIEF 25 MG  SUPP',            '0:1',  19,                 '0:2')
    ,@('BETOLVEX 1MG/ML 2 AMP',               '1:1',  64,                 '1:0')
    ,@('CYMBATEX 20 MG 30 CAPS.',             '0:2',  58,                 '0:0')
    ,@('DANSET 8MG/4ML 3 AMP.',               '0:1',  95,                 '0:3')
    ,@('DIGENORM SYRUP 120 ML',               '4:0',  55,                 '1:0')
    ,@('DOLIPRANE 1 GM 15 TABS.',             '10:0', 16,                 '0:0')
    ,@('DOSTINEX 0.5 MG 2 TABS.',             '0:0',  172,                '1:0')
    ,@('EPIMAG EFFERVESCENT 12 SACHETS',      '0:0',  26,                 '1:0')
    ,@('HEALSEC 20MG 14 CAPS',                '0:0',  47,                 '1:0')
    ,@('HERO VITAMIN DROPS',                  '1:0',  32,                 '1:0')
    ,@('INDERAL 10 MG 50 TABS',               '0:1',  45,                 '0:5')
    ,@('PANADOL ADVANCE 500 MG 48 TABLETS',   '2:0',  22.08,              '0:0')
    ,@('PANTOLOC 40MG 14 TAB',                '1:0',  51,                 '0:2')
    ,@('TRIACTIN 4MG 20 TAB',                 '1:1',  23,                 '0:2')
    ,@('جهاز محلول ',                         '10:0', 20,                 '1:0')
    ,@('سرنجات 10 سم',                        '-2:0', 8,                  '2:0')
    ,@('سرنجات 3 سم',                         '-2:0', 4,                  '2:0')
    ,@('سرنجات 5 سم',                         '-1:0', 2,                  '1:0')
    ,@('شفرات فينوس حريمي ',                  '16:0', 40,                 '2:0')
    ,@('كالونا ',                             '-1:0', 15,                 '1:0')
    ,@('كريم فيرند لافلي الصغير',             '6:0',  20,                 '1:0')
    ,@('محلول ملح',                           '27:0', 48,                 '2:0')
    ,@('مناديل سولو سحب',                     '38:0', 45,                 '1:0')
)

$firstRow = 4
$origLastRow = 23
$newLastRow = $firstRow + $finalRows.Count - 1   # 26

# --- 1) Rewrite rows 5..origLastRow in place (row 4 is unchanged) ------
for ($i = 1; $i -lt $finalRows.Count; $i++) {
    $r = $firstRow + $i
    if ($r -gt $origLastRow) { break }
    $item = $finalRows[$i]
    $ws.Range("B" + $r).Value = $item[0]
    $ws.Range("H" + $r).Value = $item[1]
    $ws.Range("L" + $r).Value = $item[2]
    $ws.Range("N" + $r).Value = $item[3]
}

# --- 2) Append brand-new rows for the items that no longer fit --------
$newRowHeights = @{ 24 = 24.75; 25 = 25.5; 26 = 24.75 }

# The totals / footer rows (originally K24:N24, A25:E25, F25:G25, I25:N25)
# get pushed down by each Insert(); re-create their merges afterwards so
# that, in the saved file, the merge records for the freshly-inserted
# rows 24-26 come *before* the (re-added) totals/footer merge records -
# matching how a from-scratch regeneration would order them.
for ($r = $origLastRow + 1; $r -le $newLastRow; $r++) {
    $ws.Rows.Item($r).Insert()
}

# The shifted totals/footer row numbers after all three inserts.
$totalsRow = $newLastRow + 1   # 27
$footerRow = $newLastRow + 2   # 28

$ws.Range("K" + $totalsRow + ":N" + $totalsRow).UnMerge()
$ws.Range("A" + $footerRow + ":E" + $footerRow).UnMerge()
$ws.Range("F" + $footerRow + ":G" + $footerRow).UnMerge()
$ws.Range("I" + $footerRow + ":N" + $footerRow).UnMerge()

for ($r = $origLastRow + 1; $r -le $newLastRow; $r++) {
    $ws.Rows.Item($r).RowHeight = $newRowHeights[$r]

    $borderRng = $ws.Range("A" + $r + ":N" + $r)
    $bottomBorder = $borderRng.Borders.Item(9)
    $bottomBorder.Color = 13882323
    $bottomBorder.LineStyle = 1

    $ws.Range("B" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()

    $idx = $r - $firstRow
    $item = $finalRows[$idx]
    $ws.Range("A" + $r).Value = $idx + 1
    $ws.Range("B" + $r).Value = $item[0]
    $ws.Range("H" + $r).Value = $item[1]
    $ws.Range("L" + $r).Value = $item[2]
    $ws.Range("N" + $r).Value = $item[3]
}

$ws.Range("K" + $totalsRow + ":N" + $totalsRow).Merge()
$ws.Range("A" + $footerRow + ":E" + $footerRow).Merge()
$ws.Range("F" + $footerRow + ":G" + $footerRow).Merge()
$ws.Range("I" + $footerRow + ":N" + $footerRow).Merge()

# --- 3) Fix up the totals row (shifted down by the 3 new rows) --------
$total = 0
foreach ($item in $finalRows) { $total += $item[2] }
$ws.Range("K" + $totalsRow).Value = $total
$ws.Rows.Item($totalsRow).RowHeight = 26.25
